# "Updated rates, added slope"
# Appends the 2022 Sep-Dec interest-rate rows (82-85) to Sheet1, continuing
# the existing year/month/y0.5/y1/y2/y3/y5 table, and moves the active
# selection to reflect the new bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# year, month, y0.5, y1, y2, y3, y5
$newRows = @(
    @(2022, "Sep", 2.47, 3.49, 3.91, 3.79, 3.75),
    @(2022, "Oct", 2.7,  3.6,  4.05, 4.01, 3.89),
    @(2022, "Nov", 2.7,  3.73, 4.08, 4.08, 3.97),
    @(2022, "Dec", 3.13, 3.8,  4.07, 4.05, 4.04)
)

$startRow = 82
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $values = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $values[0]
    $ws.Cells.Item($r, 2).Value = $values[1]
    $ws.Cells.Item($r, 3).Value = $values[2]
    $ws.Cells.Item($r, 4).Value = $values[3]
    $ws.Cells.Item($r, 5).Value = $values[4]
    $ws.Cells.Item($r, 6).Value = $values[5]
    $ws.Cells.Item($r, 7).Value = $values[6]
}

# Reflect the new bottom-of-data selection (bottom-left frozen pane).
$ws.Range("G86").Select()
